# Updates cryptos list: refresh D (price) and E (1h volume/change) columns
# for each changed row, matching the upstream scrape diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.441.82'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.581.46'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.70'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.32'
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.05'
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("E10").Value = '  -1.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0591'
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").Value = '1.807.18'
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").Value = '1.581.66'
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("E16").Value = '  -1.54%  '
$ws.Range("D17").Value = '28.454.99'
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.22'
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.67'
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("E20").Value = '  -0.62%  '
$ws.Range("D21").Value = '0.0₃0690'
$ws.Range("E21").Value = '  -2.27%  '
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.91'
$ws.Range("E23").Value = '  -3.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.14'
$ws.Range("E24").Value = '  -1.95%  '
$ws.Range("E25").Value = '  +4.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.08'
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.02'
$ws.Range("E27").Value = '  -1.49%  '
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("E29").Value = '  -2.08%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  +2.53%  '
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("E33").Value = '  -1.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.10'
$ws.Range("E34").Value = '  -2.13%  '
$ws.Range("D35").Value = '1.398.69'
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.08'
$ws.Range("E36").Value = '  +6.99%  '
$ws.Range("E37").Value = '  -4.34%  '
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.66'
$ws.Range("E39").Value = '  +2.11%  '
$ws.Range("E40").Value = '  -0.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.523'
$ws.Range("E41").Value = '  -3.66%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  +1.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.789'
$ws.Range("E44").Value = '  -2.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0465'
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("E46").Value = '  -3.12%  '
$ws.Range("E47").Value = '  -2.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.30'
$ws.Range("D49").Value = '1.718.94'
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.70'
$ws.Range("E50").Value = '  -0.69%  '
$ws.Range("D51").Value = '0.0₆0103'
$ws.Range("E51").Value = '  -0.61%  '
